$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1186.2162
$ws.Cells.Item(15, 9).Value = 1186.2162
$ws.Cells.Item(15, 11).Value = 3558.6486
$ws.Cells.Item(15, 13).Value = -3389.6486
$ws.Cells.Item(17, 8).Value = 7873.4375
$ws.Cells.Item(17, 10).Value = 7873.4375
$ws.Cells.Item(17, 12).Value = 23620.3125
$ws.Cells.Item(17, 14).Value = -23956.3125
$ws.Cells.Item(42, 8).Value = 82.333336
$ws.Cells.Item(42, 9).Value = 82.333336
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 11).Value = 247.000008
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 13).Value = -17.00000800000001
$ws.Cells.Item(42, 14).ClearContents()
$ws.Cells.Item(61, 8).Value = 4000
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 10).Value = 4000
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 12).Value = 12000
$ws.Cells.Item(61, 13).ClearContents()
$ws.Cells.Item(61, 14).Value = -12344
$ws.Cells.Item(80, 8).Value = 1254.7
$ws.Cells.Item(80, 9).Value = 986.61536
$ws.Cells.Item(80, 10).Value = 1752.5714
$ws.Cells.Item(80, 11).Value = 2959.84608
$ws.Cells.Item(80, 12).Value = 5257.7142
$ws.Cells.Item(80, 13).Value = -1961.84608
$ws.Cells.Item(80, 14).Value = -7253.7142
$ws.Cells.Item(82, 8).Value = 2775
$ws.Cells.Item(82, 9).Value = 2775
$ws.Cells.Item(82, 11).Value = 8325
$ws.Cells.Item(82, 13).Value = -7919
$ws.Cells.Item(83, 8).Value = 1254.7
$ws.Cells.Item(83, 9).Value = 986.61536
$ws.Cells.Item(83, 10).Value = 1752.5714
$ws.Cells.Item(83, 11).Value = 8879.53824
$ws.Cells.Item(83, 12).Value = 15773.1426
$ws.Cells.Item(83, 13).Value = -3887.53824
$ws.Cells.Item(83, 14).Value = -25757.1426
$ws.Cells.Item(85, 8).Value = 2775
$ws.Cells.Item(85, 9).Value = 2775
$ws.Cells.Item(85, 11).Value = 8325
$ws.Cells.Item(85, 13).Value = -6921
$ws.Cells.Item(112, 8).Value = 2215.7693
$ws.Cells.Item(112, 10).Value = 2215.7693
$ws.Cells.Item(112, 12).Value = 6647.3079
$ws.Cells.Item(112, 14).Value = -8863.3079
$ws.Cells.Item(132, 8).Value = 2177.3513
$ws.Cells.Item(132, 9).Value = 2130.3428
$ws.Cells.Item(132, 11).Value = 6391.028399999999
$ws.Cells.Item(132, 13).Value = -3861.028399999999
$ws.Cells.Item(138, 8).Value = 3864.111
$ws.Cells.Item(138, 9).Value = 1996
$ws.Cells.Item(138, 10).Value = 4500.064
$ws.Cells.Item(138, 11).Value = 5988
$ws.Cells.Item(138, 12).Value = 13500.192
$ws.Cells.Item(138, 13).Value = -848
$ws.Cells.Item(138, 14).Value = -23780.192
$ws.Cells.Item(140, 8).Value = 69998.8
$ws.Cells.Item(140, 10).Value = 69998.8
$ws.Cells.Item(140, 12).Value = 69998.8
$ws.Cells.Item(140, 14).Value = -80358.8
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5577.7163
$ws.Cells.Item(32, 9).Value = 5577.7163
$ws.Cells.Item(32, 11).Value = 5577.7163
$ws.Cells.Item(32, 13).Value = -5290.7163
$ws.Cells.Item(45, 8).Value = 2925.577
$ws.Cells.Item(45, 9).Value = 2832.0833
$ws.Cells.Item(45, 11).Value = 2832.0833
$ws.Cells.Item(45, 13).Value = -2455.0833
$ws.Cells.Item(61, 8).Value = 6898.409
$ws.Cells.Item(61, 9).Value = 6339.1665
$ws.Cells.Item(61, 11).Value = 6339.1665
$ws.Cells.Item(61, 13).Value = -6127.1665
$ws.Cells.Item(74, 8).Value = 2360.7715
$ws.Cells.Item(74, 9).Value = 1709.8
$ws.Cells.Item(74, 11).Value = 1709.8
$ws.Cells.Item(74, 13).Value = -835.8
$ws.Cells.Item(77, 8).Value = 2360.7715
$ws.Cells.Item(77, 9).Value = 1709.8
$ws.Cells.Item(77, 11).Value = 8549
$ws.Cells.Item(77, 13).Value = -4181
$ws.Cells.Item(117, 8).Value = 43333.332
$ws.Cells.Item(117, 10).Value = 43333.332
$ws.Cells.Item(117, 12).Value = 43333.332
$ws.Cells.Item(117, 14).Value = -52511.332
$ws.Cells.Item(119, 8).Value = 57183.832
$ws.Cells.Item(119, 10).Value = 57183.832
$ws.Cells.Item(119, 12).Value = 57183.832
$ws.Cells.Item(119, 14).Value = -66859.83199999999
$ws.Cells.Item(132, 8).Value = 7802.8
$ws.Cells.Item(132, 9).Value = 5000
$ws.Cells.Item(132, 10).Value = 19014
$ws.Cells.Item(132, 11).Value = 15000
$ws.Cells.Item(132, 12).Value = 57042
$ws.Cells.Item(132, 13).Value = -12470
$ws.Cells.Item(132, 14).Value = -62102
$ws.Cells.Item(136, 8).Value = 6898.409
$ws.Cells.Item(136, 9).Value = 6339.1665
$ws.Cells.Item(136, 11).Value = 19017.4995
$ws.Cells.Item(136, 13).Value = -16467.4995
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 4128.5
$ws.Cells.Item(86, 9).Value = 3359.3103
$ws.Cells.Item(86, 11).Value = 3359.3103
$ws.Cells.Item(86, 13).Value = -2236.3103
$ws.Cells.Item(89, 8).Value = 4128.5
$ws.Cells.Item(89, 9).Value = 3359.3103
$ws.Cells.Item(89, 11).Value = 16796.5515
$ws.Cells.Item(89, 13).Value = -11180.5515
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 14).ClearContents()
$ws.Cells.Item(95, 8).Value = 48749.5
$ws.Cells.Item(95, 10).Value = 48749.5
$ws.Cells.Item(95, 12).Value = 48749.5
$ws.Cells.Item(95, 14).Value = -54241.5
$ws.Cells.Item(105, 8).Value = 14110.926
$ws.Cells.Item(105, 9).Value = 19344.334
$ws.Cells.Item(105, 10).Value = 9924.200000000001
$ws.Cells.Item(105, 11).Value = 19344.334
$ws.Cells.Item(105, 12).Value = 9924.200000000001
$ws.Cells.Item(105, 13).Value = -17597.334
$ws.Cells.Item(105, 14).Value = -13418.2
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1723.55
$ws.Cells.Item(16, 9).Value = 1021.7143
$ws.Cells.Item(16, 10).Value = 3361.1667
$ws.Cells.Item(16, 11).Value = 1021.7143
$ws.Cells.Item(16, 12).Value = 3361.1667
$ws.Cells.Item(16, 13).Value = -734.7143
$ws.Cells.Item(16, 14).Value = -3935.1667
$ws.Cells.Item(94, 8).Value = 7507
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 10).Value = 7507
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 12).Value = 7507
$ws.Cells.Item(94, 13).ClearContents()
$ws.Cells.Item(94, 14).Value = -8409
$ws.Cells.Item(113, 8).Value = 1723.55
$ws.Cells.Item(113, 9).Value = 1021.7143
$ws.Cells.Item(113, 10).Value = 3361.1667
$ws.Cells.Item(113, 11).Value = 1021.7143
$ws.Cells.Item(113, 12).Value = 3361.1667
$ws.Cells.Item(113, 13).Value = 1148.2857
$ws.Cells.Item(113, 14).Value = -7701.1667
$ws.Cells.Item(122, 8).Value = 4805.1904
$ws.Cells.Item(122, 9).Value = 1715.7693
$ws.Cells.Item(122, 10).Value = 9825.5
$ws.Cells.Item(122, 11).Value = 5147.3079
$ws.Cells.Item(122, 12).Value = 29476.5
$ws.Cells.Item(122, 13).Value = -2697.3079
$ws.Cells.Item(122, 14).Value = -34376.5
$ws.Cells.Item(133, 8).Value = 56818.832
$ws.Cells.Item(133, 10).Value = 57425.11
$ws.Cells.Item(133, 12).Value = 57425.11
$ws.Cells.Item(133, 14).Value = -62485.11
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 235.28572
$ws.Cells.Item(12, 10).Value = 336
$ws.Cells.Item(12, 12).Value = 1008
$ws.Cells.Item(12, 14).Value = -1354
$ws.Cells.Item(113, 8).Value = 2368.25
$ws.Cells.Item(113, 9).Value = 1988.6
$ws.Cells.Item(113, 10).Value = 3001
$ws.Cells.Item(113, 11).Value = 5965.799999999999
$ws.Cells.Item(113, 12).Value = 9003
$ws.Cells.Item(113, 13).Value = -3795.799999999999
$ws.Cells.Item(113, 14).Value = -13343
$ws.Cells.Item(124, 8).Value = 4610.25
$ws.Cells.Item(124, 10).Value = 5776.6
$ws.Cells.Item(124, 12).Value = 17329.8
$ws.Cells.Item(124, 14).Value = -27149.8
$ws.Cells.Item(128, 8).Value = 202580.4
$ws.Cells.Item(128, 9).Value = 202580.4
$ws.Cells.Item(128, 11).Value = 607741.2
$ws.Cells.Item(128, 13).Value = -602761.2
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(101, 8).Value = 5657
$ws.Cells.Item(101, 10).Value = 5657
$ws.Cells.Item(101, 12).Value = 5657
$ws.Cells.Item(101, 14).Value = -12147
$ws.Cells.Item(102, 8).Value = 5764.4614
$ws.Cells.Item(102, 9).Value = 1820.8334
$ws.Cells.Item(102, 10).Value = 9144.714
$ws.Cells.Item(102, 11).Value = 1820.8334
$ws.Cells.Item(102, 12).Value = 9144.714
$ws.Cells.Item(102, 13).Value = -198.8334
$ws.Cells.Item(102, 14).Value = -12388.714
$ws.Cells.Item(141, 8).Value = 63853.555
$ws.Cells.Item(141, 10).Value = 63853.555
$ws.Cells.Item(141, 12).Value = 63853.555
$ws.Cells.Item(141, 14).Value = -74213.55499999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 3082.353
$ws.Cells.Item(93, 9).Value = 3135.7144
$ws.Cells.Item(93, 10).Value = 2833.3333
$ws.Cells.Item(93, 11).Value = 3135.7144
$ws.Cells.Item(93, 12).Value = 2833.3333
$ws.Cells.Item(93, 13).Value = -1887.7144
$ws.Cells.Item(93, 14).Value = -5329.3333
$ws.Cells.Item(123, 8).Value = 68000
$ws.Cells.Item(123, 10).Value = 68000
$ws.Cells.Item(123, 12).Value = 68000
$ws.Cells.Item(123, 14).Value = -77800
$ws.Cells.Item(136, 8).Value = 5823.385
$ws.Cells.Item(136, 9).Value = 3671.4285
$ws.Cells.Item(136, 11).Value = 11014.2855
$ws.Cells.Item(136, 13).Value = -8464.2855
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 4847.4614
$ws.Cells.Item(132, 9).Value = 3101.5
$ws.Cells.Item(132, 10).Value = 7641
$ws.Cells.Item(132, 11).Value = 9304.5
$ws.Cells.Item(132, 12).Value = 22923
$ws.Cells.Item(132, 13).Value = -6774.5
$ws.Cells.Item(132, 14).Value = -27983
$ws.Cells.Item(141, 8).Value = 196904.33
$ws.Cells.Item(141, 10).Value = 196904.33
$ws.Cells.Item(141, 12).Value = 196904.33
$ws.Cells.Item(141, 14).Value = -207264.33
